$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Old -> New expression mapping, in the exact order the expressions
# appear in the table (row-major, skipping blank spacer rows).
$replacements = @(
    "69÷3=", "91÷5=",
    "92÷4=", "86÷6=",
    "77÷9=", "14÷2=",
    "62÷6=", "37÷5=",
    "28÷2=", "19÷5=",
    "59÷8=", "94÷4=",
    "41÷5=", "17÷7=",
    "52÷3=", "27÷2=",
    "27÷3=", "16÷8=",
    "69÷7=", "10÷2=",
    "32÷6=", "67÷6=",
    "33÷8=", "60÷5=",
    "43÷6=", "82÷7=",
    "33÷2=", "20÷5=",
    "65÷9=", "81÷7=",
    "87÷4=", "98÷7=",
    "63÷4=", "77÷4=",
    "69÷6=", "38÷3=",
    "93÷8=", "71÷2=",
    "65÷7=", "16÷5=",
    "25÷8=", "46÷9=",
    "81÷9=", "76÷2=",
    "86÷6=", "63÷5=",
    "48÷2=", "49÷8=",
    "97÷7=", "54÷4="
)

# Find.Execute operates on the whole document rather than being confined
# to a cell's Range, which breaks naive search/replace when an output
# value collides with another cell's input value (e.g. "86÷6=" is both
# a source and a target here). So instead assign each cell's Range.Text
# directly -- that mutation really is scoped to the individual cell.
$idx = 0
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    for ($j = 1; $j -le $tbl.Columns.Count; $j++) {
        $cell = $tbl.Cell($i, $j)
        $txt = $cell.Range.Text
        if ($txt -match "^\d+÷\d+=") {
            $old = $replacements[$idx]
            $new = $replacements[$idx + 1]
            if ($txt.Substring(0, $old.Length) -eq $old) {
                $cell.Range.Text = $new
            }
            $idx += 2
        }
    }
}
